# "Corrections made for review 3"
#
# 1. Relabel the T-critical row as a two-tailed value.
# 2. Fix the sign of the per-row difference formula in C2 (it was using the
#    shared-formula master B2-A2; the corrected version is A2-B2). Every
#    dependent statistic below (mean of differences, std dev, SE, t-stat)
#    recalculates automatically from this single change.
# 3. Update the hard-coded T-critical lookup value to the two-tailed
#    critical value (df = 23, alpha = .05).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Inferential Stat. Analysis")

# Relabel "Tcritical" -> "Tcritical (2-tailed)"
$ws.Range("E21").Value = "Tcritical (2-tailed)"

# Correct the C2 difference formula: was =B2-A2, now =A2-B2
$ws.Range("C2").Formula = "=A2-B2"

# Update the T-critical table lookup value for the two-tailed test
$ws.Range("E22").Value = 2.069
